$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D for the new fiscal quarters (2018-12-29 and 2018-09-29)
$ws.Columns("D:E").Insert()

# Populate the two new quarter columns with their values for every data row
$ws.Range("D7").Value = 43463
$ws.Range("E7").Value = 43372
$ws.Range("D8").Value = 203200
$ws.Range("E8").Value = 197900
$ws.Range("D9").Value = 124800
$ws.Range("E9").Value = 120600
$ws.Range("D10").Value = 78400
$ws.Range("E10").Value = 77300
$ws.Range("D12").Value = 7200
$ws.Range("E12").Value = 8200
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = 1400
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 185200
$ws.Range("E17").Value = 182100
$ws.Range("D18").Value = 18000
$ws.Range("E18").Value = 15800
$ws.Range("D20").Value = -6800
$ws.Range("E20").Value = -1300
$ws.Range("D21").Value = 20200
$ws.Range("E21").Value = 23200
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 11200
$ws.Range("E23").Value = 14500
$ws.Range("D24").Value = 2000
$ws.Range("E24").Value = 3400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 9200
$ws.Range("E26").Value = 11100
$ws.Range("D27").Value = 9200
$ws.Range("E27").Value = 11100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 1300
$ws.Range("E29").Value = -400
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 6800
$ws.Range("E32").Value = 1300
$ws.Range("D33").Value = 10500
$ws.Range("E33").Value = 10800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 10500
$ws.Range("E35").Value = 10800
$ws.Range("D38").Value = 43463
$ws.Range("E38").Value = 43372
$ws.Range("D41").Value = 70400
$ws.Range("E41").Value = 71800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 176000
$ws.Range("E43").Value = 192700
$ws.Range("D44").Value = 183500
$ws.Range("E44").Value = 139100
$ws.Range("D45").Value = 30100
$ws.Range("E45").Value = 24600
$ws.Range("D46").Value = 460000
$ws.Range("E46").Value = 428200
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 89300
$ws.Range("E48").Value = 90300
$ws.Range("D49").Value = 690200
$ws.Range("E49").Value = 615400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 8300
$ws.Range("E52").Value = 5500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1247900
$ws.Range("E54").Value = 1139400
$ws.Range("D57").Value = 49200
$ws.Range("E57").Value = 47900
$ws.Range("D58").Value = 30500
$ws.Range("E58").Value = 32700
$ws.Range("D59").Value = 177600
$ws.Range("E59").Value = 158500
$ws.Range("D60").Value = 257300
$ws.Range("E60").Value = 239100
$ws.Range("D61").Value = 435600
$ws.Range("E61").Value = 355600
$ws.Range("D62").Value = 81700
$ws.Range("E62").Value = 66700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 774500
$ws.Range("E66").Value = 661500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 299500
$ws.Range("E72").Value = 300600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 473300
$ws.Range("E76").Value = 477900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43463
$ws.Range("E80").Value = 43372
$ws.Range("D81").Value = 10500
$ws.Range("E81").Value = 10800
$ws.Range("D83").Value = 9000
$ws.Range("E83").Value = 8600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 10600
$ws.Range("E89").Value = 11100
$ws.Range("D91").Value = -3800
$ws.Range("E91").Value = -2500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -82100
$ws.Range("E94").Value = 4200
$ws.Range("D96").Value = -5400
$ws.Range("E96").Value = -5400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 70600
$ws.Range("E100").Value = -9000
$ws.Range("D101").Value = -500
$ws.Range("E101").Value = -1000
$ws.Range("D102").Value = -1400
$ws.Range("E102").Value = 5400

# Copy number formats from column F (the old column D, now shifted) onto the new D:E columns
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data corrections carried over from the old sheet (values updated at their new shifted positions)
$ws.Range("I89").Value = 19200
$ws.Range("I91").Value = -4600
$ws.Range("J91").Value = -3900
